# Reverted block-based CS. Created LS vs CS interferogram.
#
# - Add a new "2dmock1.fits/.csv" row (row 5) to the "2D" sheet's table,
#   mirroring the 1D mocks' errors vs samples data.
# - Switch the active/selected tab from "1D NEW" to "2D".

$wb = $excel.ActiveWorkbook

# "2D" sheet is the 4th tab (1D OLD, 1D NEW, 1D TRAINING, 2D).
$ws2d = $wb.Worksheets.Item(4)

# Append the new data row under the existing "2dmock.fits" row.
$ws2d.Range("B5").Value = "2dmock1.fits/.csv"
$ws2d.Range("C5").Value = 100
$ws2d.Range("D5").Value = 75
$ws2d.Range("E5").Value = 200
$ws2d.Range("F5").Value = 3
$ws2d.Range("G5").Value = 3
$ws2d.Range("H5").Value = 30
$ws2d.Range("I5").Value = 0
$ws2d.Range("J5").Value = 20
$ws2d.Range("K5").Value = 30

# Make "2D" the active sheet/tab, with K6 selected.
[void]$ws2d.Activate()
[void]$ws2d.Range("K6").Select()
